$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Web")

# Row 19: UM_searchUserRole
$ws.Range("A19").Value = "UM_searchUserRole"
$ws.Range("B19").Value = "//*[@id=""app""]/div[1]/div[2]/div[2]/div/div[1]/div[2]/form/div[1]/div/div[2]/div/div[2]/div/div"
$ws.Range("C19").Value = "By.xpath"

# Row 20: UM_selectUserRole (re-uses same xpath as row 19)
$ws.Range("A20").Value = "UM_selectUserRole"
$ws.Range("B20").Value = "//*[@id=""app""]/div[1]/div[2]/div[2]/div/div[1]/div[2]/form/div[1]/div/div[2]/div/div[2]/div/div"
$ws.Range("C20").Value = "By.xpath"

# Row 21: UM_searchResultUserRole
$ws.Range("A21").Value = "UM_searchResultUserRole"
$ws.Range("B21").Value = "//*[@id=""app""]/div[1]/div[2]/div[2]/div/div[2]/div[3]/div/div[2]/div[1]/div/div[3]/div"
$ws.Range("C21").Value = "By.xpath"
